$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 291 (shifts old rows 291-304 down to 293-306)
$ws.Rows("291:292").Insert()

# New row 291: Pera, Packham's Triumph, Primera - week of 2021-11-09 (serial 44509)
$ws.Range("A291").Value = 11
$ws.Range("B291").Value = "Vega Monumental Concepción"
$ws.Range("C291").Value = "Bíobío"
$ws.Range("D291").Value = 44509
$ws.Range("E291").Value = 8
$ws.Range("F291").Value = "Fruta"
$ws.Range("G291").Value = 100104
$ws.Range("H291").Value = "Frutos de pepita"
$ws.Range("I291").Value = 100104005
$ws.Range("J291").Value = "Pera"
$ws.Range("K291").Value = "Packham's Triumph"
$ws.Range("L291").Value = "Primera"
$ws.Range("M291").Value = 250
$ws.Range("N291").Value = 10000
$ws.Range("O291").Value = 10000
$ws.Range("P291").Value = 10000
$ws.Range("Q291").Value = "`$/caja 17 kilos empedrada"
$ws.Range("R291").Value = "Región de O'Higgins"
$ws.Range("S291").Value = 588
$ws.Range("T291").Value = 17

# New row 292: Pera, Packham's Triumph, Segunda - week of 2021-11-09 (serial 44509)
$ws.Range("A292").Value = 11
$ws.Range("B292").Value = "Vega Monumental Concepción"
$ws.Range("C292").Value = "Bíobío"
$ws.Range("D292").Value = 44509
$ws.Range("E292").Value = 8
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100104
$ws.Range("H292").Value = "Frutos de pepita"
$ws.Range("I292").Value = 100104005
$ws.Range("J292").Value = "Pera"
$ws.Range("K292").Value = "Packham's Triumph"
$ws.Range("L292").Value = "Segunda"
$ws.Range("M292").Value = 250
$ws.Range("N292").Value = 8000
$ws.Range("O292").Value = 8000
$ws.Range("P292").Value = 8000
$ws.Range("Q292").Value = "`$/caja 17 kilos empedrada"
$ws.Range("R292").Value = "Región de O'Higgins"
$ws.Range("S292").Value = 471
$ws.Range("T292").Value = 17
